$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column widths for Z (26) and AA (27) ---------------------------
$ws.Columns.Item(26).ColumnWidth = 50.33
$ws.Columns.Item(27).ColumnWidth = 41.33

# --- New header / value cells in columns Z and AA ------------------------
$ws.Range("Z2").Value = "Number of months since Policy issued"
$ws.Range("AA2").Value = "Payment due to Death"
$ws.Range("Z3").Value = "Death Benefit Payment in Percentage"
$ws.Range("Z4").Value = "Death Method used for Death Benefit Payout"

# --- Formatting: thin border, vertically centered + wrapped text, ------
# --- slightly smaller font -------------------------------------------
$anchor = $ws.Range("Z2")
$anchor.Font.Size = 7.5
$anchor.Borders.LineStyle = 1
$anchor.VerticalAlignment = -4108
$anchor.WrapText = $true

# Propagate the exact same formatting to the other new cells by copying
# the already-built format (keeps the style table minimal/clean).
$anchor.Copy() | Out-Null
$ws.Range("AA2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Selection / view state, matching where the author ended up --------
$ws.Range("Z2:AA4").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 14 | Out-Null
